$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold plain text (not real numbers - e.g. thousands
# separators like "29.745.18" or padded percentage strings). Force the
# affected cells to be treated as text first so Excel doesn't auto-convert
# numeric-looking strings (like "0.9988") into floating point numbers, then
# strip the temporary number format back off again once the values are set.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.745.18"
$ws.Range("E2").Value = "  +5.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.921.98"
$ws.Range("E3").Value = "  +3.55%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.29%  "

# Row 5 - BNB
$ws.Range("D5").Value = "335.39"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6 - USDC
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.21%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +3.10%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.4119"
$ws.Range("E8").Value = "  +5.17%  "

# Row 9 - OKB
$ws.Range("D9").Value = "48.13"
$ws.Range("E9").Value = "  +1.63%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.08035"
$ws.Range("E10").Value = "  +3.32%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  +3.52%  "

# Row 12 - Solana
$ws.Range("D12").Value = "22.48"
$ws.Range("E12").Value = "  +5.66%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.930.32"
$ws.Range("E13").Value = "  +3.64%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.996"
$ws.Range("E14").Value = "  +3.89%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.190"
$ws.Range("E15").Value = "  +3.71%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "90.13"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +2.35%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  +1.03%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +5.17%  "

# Row 21 - Dai
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "29.697.37"
$ws.Range("E22").Value = "  +5.07%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.581"
$ws.Range("E23").Value = "  +5.90%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "11.76"
$ws.Range("E24").Value = "  +10.60%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.207"
$ws.Range("E25").Value = "  -1.90%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.168.37"
$ws.Range("E26").Value = "  +4.52%  "

# Row 27 - Monero
$ws.Range("D27").Value = "156.22"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  +4.10%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  +5.47%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "5.725"
$ws.Range("E30").Value = "  +8.93%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "117.60"
$ws.Range("E31").Value = "  +1.33%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +14.44%  "

# Row 33 - Stellar
$ws.Range("D33").Value = "0.09475"
$ws.Range("E33").Value = "  +2.44%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.441"
$ws.Range("E34").Value = "  +5.19%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "3.568"
$ws.Range("E35").Value = "  -1.07%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  +4.90%  "

# Row 37 - was Hedera, now VeChain
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02276"
$ws.Range("E37").Value = "  +4.01%  "

# Row 38 - was VeChain, now Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06141"
$ws.Range("E38").Value = "  +2.53%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "8.457"
$ws.Range("E39").Value = "  +3.52%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "1.182"

# Row 41 - TheSandbox
$ws.Range("D41").Value = "0.5897"
$ws.Range("E41").Value = "  +4.54%  "

# Row 42 - Algorand
$ws.Range("D42").Value = "0.1849"
$ws.Range("E42").Value = "  +3.67%  "

# Row 43 - Aptos
$ws.Range("D43").Value = "10.26"
$ws.Range("E43").Value = "  +3.27%  "

# Row 44 - WEMIXTOKEN
$ws.Range("D44").Value = "1.260"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45 - RenderToken
$ws.Range("D45").Value = "2.363"
$ws.Range("E45").Value = "  +2.89%  "

# Row 46 - Cronos
$ws.Range("D46").Value = "0.07512"
$ws.Range("E46").Value = "  +4.71%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.5583"
$ws.Range("E47").Value = "  +4.38%  "

# Row 48 - EnergySwap
$ws.Range("E48").Value = "  +4.03%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  +4.10%  "

# Row 50 - Quant
$ws.Range("D50").Value = "113.32"
$ws.Range("E50").Value = "  +3.77%  "

# Row 51 - WOONetwork
$ws.Range("D51").Value = "0.3017"
$ws.Range("E51").Value = "  +15.31%  "

# Remove the temporary text formatting again so styles are left untouched.
$textRange.ClearFormats()
